# Trade #18 closed at 2026-02-16 21:25:06 - leadlag DOWN +0.000%
# Append a new trade row to the bottom of the "leadlag" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 17

# Trade #
$ws.Cells.Item($row, 1).Value = 18

# Date - stored as plain text in the sheet, guard against Excel's
# automatic date inference by using the quote-prefix trick, then
# strip the resulting cell style back to Normal so no extra
# formatting is left behind.
$ws.Cells.Item($row, 2).Value = "'2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"

# Time
$ws.Cells.Item($row, 3).Value = "21:25:06"

# Strategy
$ws.Cells.Item($row, 4).Value = "leadlag"

# Side
$ws.Cells.Item($row, 5).Value = "DOWN"

# Entry Price
$ws.Cells.Item($row, 6).Value = 69213.42

# Exit Price (blank - trade still open)
$ws.Cells.Item($row, 7).Value = "'"
$ws.Cells.Item($row, 7).Style = "Normal"

# Status
$ws.Cells.Item($row, 8).Value = "OPEN"

# P&L %
$ws.Cells.Item($row, 9).Value = 0

# P&L $
$ws.Cells.Item($row, 10).Value = 0

# Confidence
$ws.Cells.Item($row, 11).Value = 0.75

# Entry Reason
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.103% move"

# Exit Reason (blank - trade still open)
$ws.Cells.Item($row, 13).Value = "'"
$ws.Cells.Item($row, 13).Style = "Normal"

# Duration (min)
$ws.Cells.Item($row, 14).Value = 0
